$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers, matching column A format)
$newRows = @(
    @{ Row = 230; A = 44304; B = 4; C = 26; D = 172.3337973089415 },
    @{ Row = 231; A = 44305; B = 1; C = 17; D = 112.679790548154 },
    @{ Row = 232; A = 44306; B = 0; C = 16; D = 106.0515675747332 },
    @{ Row = 233; A = 44307; B = 0; C = 15; D = 99.42334460131239 }
)

# Reuse formatting from the last existing data row (229) so the new cells
# keep the same style (date format, borders, alignment) as the rest of
# column A, instead of creating a brand new style entry.
$lastDataRow = 229

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Range("A$lastDataRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
}

$excel.CutCopyMode = 0
